$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# --- Locate existing sheets ---
$total = $wb.Worksheets.Item(1)     # "总计"
$q2    = $wb.Worksheets.Item(2)     # "2022-Q2" (will become the new "2022-Q3" sheet)

# Rename first (data/formatting untouched), then duplicate so the exact
# original "2022-Q2" data/styles survive unchanged as a brand new sheet
# placed right after it.
$q2.Name = "2022-Q3"
$q2.Copy($null, $q2)
$q2Copy = $wb.Worksheets.Item(3)
$q2Copy.Name = "2022-Q2"

# --- Overwrite the (renamed) sheet's data with the Q3 fund-holdings figures ---

# Header row formatting matches the "总计" sheet's header style.
$total.Range("B1:D1").Copy()
$q2.Range("B1:H1").PasteSpecial($xlPasteFormats)

$q2.Range("B1").Value = "基金代码"
$q2.Range("C1").Value = "基金名称"
$q2.Range("D1").Value = "基金规模"
$q2.Range("E1").Value = "股票总仓位"
$q2.Range("F1").Value = "仓位占比"
$q2.Range("G1").Value = "持有市值(亿元)"
$q2.Range("H1").Value = "仓位排名"

# Index column (A2:A4) matches the "总计" sheet's index-cell style.
$total.Range("A2").Copy()
$q2.Range("A2:A4").PasteSpecial($xlPasteFormats)

# Fund-code / decimal-looking columns must stay text (preserve leading
# zeros, e.g. "011903", and the literal decimal strings like "1.56").
$q2.Range("B2:G4").NumberFormat = "@"

$q2.Range("A2").Value = 0
$q2.Range("B2").Value = "011903"
$q2.Range("C2").Value = "南方领航优选混合A"
$q2.Range("D2").Value = "1.56"
$q2.Range("E2").Value = "82.01"
$q2.Range("F2").Value = "5.61"
$q2.Range("G2").Value = "0.0875"
$q2.Range("H2").Value = 4

$q2.Range("A3").Value = 1
$q2.Range("B3").Value = "004703"
$q2.Range("C3").Value = "南方兴盛先锋灵活配置混合"
$q2.Range("D3").Value = "0.97"
$q2.Range("E3").Value = "82.88"
$q2.Range("F3").Value = "5.43"
$q2.Range("G3").Value = "0.0527"
$q2.Range("H3").Value = 5

$q2.Range("A4").Value = 2
$q2.Range("B4").Value = "011904"
$q2.Range("C4").Value = "南方领航优选混合C"
$q2.Range("D4").Value = "0.48"
$q2.Range("E4").Value = "82.01"
$q2.Range("F4").Value = "5.61"
$q2.Range("G4").Value = "0.0269"
$q2.Range("H4").Value = 4

# --- Update the "总计" summary sheet: existing row becomes the Q3 total,
#     and a new row is appended for the (now second) Q2 total. ---
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.17

$total.Range("A2").Copy($total.Range("A3"))
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.05

# Keep "总计" as the active/visible sheet, matching the original workbook state.
$total.Activate()
